# Applies the cryptos list update committed on Thu Oct 26 16:39:15 UTC 2023
# Updates Price (column D) and Volume(1h) (column E) for rows 2-51,
# and swaps the FraxShare/Kaspa rows (44/45) to reflect their new rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'33.951.81"
$ws.Range("E2").Value = "'  -2.29%  "
$ws.Range("D3").Value = "'1.771.01"
$ws.Range("E3").Value = "'  -1.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'220.90"
$ws.Range("E5").Value = "'  -1.95%  "
$ws.Range("E6").Value = "'  -1.64%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  -0.21%  "
$ws.Range("D8").Value = "'30.94"
$ws.Range("E8").Value = "'  -6.43%  "
$ws.Range("E9").Value = "'  -0.34%  "
$ws.Range("E10").Value = "'  +5.06%  "
$ws.Range("E11").Value = "'  -1.76%  "
$ws.Range("D12").Value = "'2.026.78"
$ws.Range("E12").Value = "'  -1.43%  "
$ws.Range("D13").Value = "'1.769.76"
$ws.Range("E13").Value = "'  -1.29%  "
$ws.Range("D14").Value = "'10.49"
$ws.Range("E14").Value = "'  -5.44%  "
$ws.Range("E15").Value = "'  -2.49%  "
$ws.Range("D16").Value = "'33.936.88"
$ws.Range("E16").Value = "'  -2.28%  "
$ws.Range("D17").Value = "'4.19"
$ws.Range("E17").Value = "'  -2.68%  "
$ws.Range("D18").Value = "'67.59"
$ws.Range("E18").Value = "'  -2.60%  "
$ws.Range("D19").Value = "'243.03"
$ws.Range("E19").Value = "'  -5.33%  "
$ws.Range("D20").Value = "'0.0₃0774"
$ws.Range("E20").Value = "'  +1.31%  "
$ws.Range("E21").Value = "'  -0.04%  "
$ws.Range("E22").Value = "'  +0.78%  "
$ws.Range("D23").Value = "'4.01"
$ws.Range("E23").Value = "'  -5.25%  "
$ws.Range("E24").Value = "'  -1.28%  "
$ws.Range("D25").Value = "'157.51"
$ws.Range("E25").Value = "'  -0.91%  "
$ws.Range("D26").Value = "'16.31"
$ws.Range("E26").Value = "'  -0.88%  "
$ws.Range("D27").Value = "'6.97"
$ws.Range("E27").Value = "'  -2.20%  "
$ws.Range("E28").Value = "'  -2.47%  "
$ws.Range("E29").Value = "'  -0.09%  "
$ws.Range("D30").Value = "'0.0518"
$ws.Range("E30").Value = "'  -0.43%  "
$ws.Range("D31").Value = "'3.70"
$ws.Range("E31").Value = "'  -2.01%  "
$ws.Range("E32").Value = "'  +0.02%  "
$ws.Range("D33").Value = "'3.50"
$ws.Range("E33").Value = "'  -3.01%  "
$ws.Range("E34").Value = "'  -3.55%  "
$ws.Range("D35").Value = "'1.395.72"
$ws.Range("E35").Value = "'  -4.18%  "
$ws.Range("E36").Value = "'  -1.60%  "
$ws.Range("D37").Value = "'0.631"
$ws.Range("E37").Value = "'  +0.08%  "
$ws.Range("E38").Value = "'  -2.47%  "
$ws.Range("D39").Value = "'0.927"
$ws.Range("E39").Value = "'  +2.75%  "
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("D41").Value = "'78.60"
$ws.Range("E41").Value = "'  -6.19%  "
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "'  -5.26%  "
$ws.Range("D43").Value = "'2.10"
$ws.Range("E43").Value = "'  +0.82%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.87"
$ws.Range("E44").Value = "'  -1.26%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.0488"
$ws.Range("E45").Value = "'  -3.65%  "
$ws.Range("D46").Value = "'1.03"
$ws.Range("E46").Value = "'  -1.18%  "
$ws.Range("D47").Value = "'1.924.31"
$ws.Range("E47").Value = "'  -1.66%  "
$ws.Range("D48").Value = "'104.09"
$ws.Range("E48").Value = "'  -1.01%  "
$ws.Range("E49").Value = "'  -0.57%  "
$ws.Range("D50").Value = "'11.82"
$ws.Range("E50").Value = "'  -0.89%  "
$ws.Range("E51").Value = "'  -1.75%  "
